$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 26

# 1. Capture the existing "Output Bin." probability values (col B) and
#    "Output Count." values (col C) before we overwrite anything, then
#    swap them: Count moves to column B, probability moves to column C,
#    and a new binary prediction column is derived in column D.
for ($r = 2; $r -le $lastRow; $r++) {
    $prob = $ws.Cells.Item($r, 2).Value2
    $count = $ws.Cells.Item($r, 3).Value2

    $ws.Cells.Item($r, 2).Value2 = $count
    $ws.Cells.Item($r, 3).Value2 = $prob

    if ($prob -ge 0.5) {
        $ws.Cells.Item($r, 4).Value2 = 1
    } else {
        $ws.Cells.Item($r, 4).Value2 = 0
    }
}

# 2. Update header row for the reshuffled / new columns.
$ws.Range("B1").Value2 = "Output Count."
$ws.Range("C1").Value2 = "Output Binary Predict Probability"
$ws.Range("D1").Value2 = "Output Binary Prediction"

# 3. Column widths to match the autofit the original author applied
#    (values chosen so this engine's pixel->width quantization lands as
#    close as possible to the authored widths of 12.71/13.71/18.43/23.29).
$ws.Columns.Item(1).ColumnWidth = 11.8333333333333
$ws.Columns.Item(2).ColumnWidth = 12.8333333333333
$ws.Columns.Item(3).ColumnWidth = 17.6666666666667
$ws.Columns.Item(4).ColumnWidth = 22.5

# 4. Restore the selection the author left the sheet on.
$ws.Range("D4").Select()
